$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 11.319211
$ws.Range("H2").Value = 33.957633
$ws.Range("I2").Value = 0.09922284194232082
$ws.Range("J2").Value = 0.09922284194232082
$ws.Range("M2").Value = 55.60163533333334
$ws.Range("N2").Value = 166.804906
$ws.Range("O2").Value = 0.2980784899567677
$ws.Range("P2").Value = 0.2980784899567676
$ws.Range("Q2").Value = 629.3666422830555
$ws.Range("R2").Value = 5664.299780547499
$ws.Range("S2").Value = 0.02957619489538602
$ws.Range("T2").Value = 0.02957619489538602
# Row 3
$ws.Range("G3").Value = 11.319211
$ws.Range("H3").Value = 33.957633
$ws.Range("I3").Value = 0.09922284194232082
$ws.Range("J3").Value = 0.09922284194232082
$ws.Range("O3").Value = 0.5359390331190738
$ws.Range("P3").Value = 0.5359390331190738
$ws.Range("Q3").Value = 1131.588360473444
$ws.Range("R3").Value = 10184.295244261
$ws.Range("S3").Value = 0.0531773939738941
$ws.Range("T3").Value = 0.0531773939738941
# Row 4
$ws.Range("G4").Value = 11.319211
$ws.Range("H4").Value = 33.957633
$ws.Range("I4").Value = 0.09922284194232082
$ws.Range("J4").Value = 0.09922284194232082
$ws.Range("M4").Value = 30.96129866666666
$ws.Range("N4").Value = 92.88389599999999
$ws.Range("O4").Value = 0.1659824769241586
$ws.Range("P4").Value = 0.1659824769241586
$ws.Range("Q4").Value = 350.4574724420187
$ws.Range("R4").Value = 3154.117251978168
$ws.Range("S4").Value = 0.0164692530730407
$ws.Range("T4").Value = 0.0164692530730407
# Row 5
$ws.Range("I5").Value = 0.3843080175847637
$ws.Range("J5").Value = 0.3843080175847637
$ws.Range("M5").Value = 55.60163533333334
$ws.Range("N5").Value = 166.804906
$ws.Range("O5").Value = 0.2980784899567677
$ws.Range("P5").Value = 0.2980784899567676
$ws.Range("Q5").Value = 2437.650866424304
$ws.Range("R5").Value = 21938.85779781874
$ws.Range("S5").Value = 0.1145539535599453
$ws.Range("T5").Value = 0.1145539535599453
# Row 6
$ws.Range("I6").Value = 0.3843080175847637
$ws.Range("J6").Value = 0.3843080175847637
$ws.Range("O6").Value = 0.5359390331190738
$ws.Range("P6").Value = 0.5359390331190738
$ws.Range("S6").Value = 0.2059656673642862
$ws.Range("T6").Value = 0.2059656673642862
# Row 7
$ws.Range("I7").Value = 0.3843080175847637
$ws.Range("J7").Value = 0.3843080175847637
$ws.Range("M7").Value = 30.96129866666666
$ws.Range("N7").Value = 92.88389599999999
$ws.Range("O7").Value = 0.1659824769241586
$ws.Range("P7").Value = 0.1659824769241586
$ws.Range("S7").Value = 0.06378839666053218
$ws.Range("T7").Value = 0.06378839666053218
# Row 8
$ws.Range("I8").Value = 0.5164691404729155
$ws.Range("J8").Value = 0.5164691404729155
$ws.Range("M8").Value = 55.60163533333334
$ws.Range("N8").Value = 166.804906
$ws.Range("O8").Value = 0.2980784899567677
$ws.Range("P8").Value = 0.2980784899567676
$ws.Range("Q8").Value = 3275.943748630061
$ws.Range("R8").Value = 29483.49373767055
$ws.Range("S8").Value = 0.1539483415014364
$ws.Range("T8").Value = 0.1539483415014363
# Row 9
$ws.Range("I9").Value = 0.5164691404729155
$ws.Range("J9").Value = 0.5164691404729155
$ws.Range("O9").Value = 0.5359390331190738
$ws.Range("P9").Value = 0.5359390331190738
$ws.Range("S9").Value = 0.2767959717808934
$ws.Range("T9").Value = 0.2767959717808934
# Row 10
$ws.Range("I10").Value = 0.5164691404729155
$ws.Range("J10").Value = 0.5164691404729155
$ws.Range("M10").Value = 30.96129866666666
$ws.Range("N10").Value = 92.88389599999999
$ws.Range("O10").Value = 0.1659824769241586
$ws.Range("P10").Value = 0.1659824769241586
$ws.Range("Q10").Value = 1824.181468916775
$ws.Range("S10").Value = 0.08572482719058572
$ws.Range("T10").Value = 0.08572482719058572
